# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
# Both sheets hold the same set of exhibition rows (mirrored data), so the
# same F-column values need to be bumped on each, though the row numbers
# differ slightly between the two sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14842
$wsExhibit.Range("F3").Value = 18416
$wsExhibit.Range("F14").Value = 98
$wsExhibit.Range("F22").Value = 7631
$wsExhibit.Range("F28").Value = 5942
$wsExhibit.Range("F34").Value = 5278
$wsExhibit.Range("F35").Value = 25

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14842
$wsAll.Range("F3").Value = 18416
$wsAll.Range("F14").Value = 98
$wsAll.Range("F23").Value = 7631
$wsAll.Range("F30").Value = 5942
$wsAll.Range("F36").Value = 5278
$wsAll.Range("F37").Value = 25
